$d = $word.ActiveDocument

# The cover page has a standalone paragraph whose entire text is "202"
# (right after "UNIVERSITAS MATARAM"). Another, unrelated heading already
# contains "...STRUKTUR DATA 2024", so a plain substring Find would risk
# matching the wrong spot. Using MatchWholeWord restricts the search to the
# isolated "202" token and finds exactly the cover-page run.
$rng = $d.Content
$found = $rng.Find.Execute("202", $false, $true, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the end of the found "202" and type the missing "4",
    # turning the year into "2024" - same as a user clicking right after
    # "202" and typing the last digit.
    $rng.Collapse(0)
    $rng.Select()
    $word.Selection.TypeText("4")
}
